$wb = $excel.ActiveWorkbook

# --- Sheet 1: "species mean" (A1:H3) ---
$ws1 = $wb.Worksheets.Item("species mean")
$ws1.Range("E2").Value = 0.5132529441592365
$ws1.Range("F2").Value = 9.567735170012789
$ws1.Range("G2").Value = 0.1915974366617081
$ws1.Range("H2").Value = 499.3890768034377
$ws1.Range("E3").Value = 1.044385408163225
$ws1.Range("F3").Value = 30.33532670799015
$ws1.Range("G3").Value = 0.1874338713503914
$ws1.Range("H3").Value = 69.81732006839894

# --- Sheet 2: "species bioshaker mean" (A1:I5) ---
$ws2 = $wb.Worksheets.Item("species bioshaker mean")
$ws2.Range("F2").Value = 0.6281494734143215
$ws2.Range("G2").Value = 7.436949464255736
$ws2.Range("H2").Value = 0.2145031977514762
$ws2.Range("I2").Value = 667.0968230247564
$ws2.Range("F3").Value = 0.3966415413331798
$ws2.Range("G3").Value = 11.73032364749756
$ws2.Range("H3").Value = 0.1683497985407496
$ws2.Range("I3").Value = 329.1782298922487
$ws2.Range("F4").Value = 0.3597759374553952
$ws2.Range("G4").Value = 18.01459785928747
$ws2.Range("H4").Value = 0.2296802866683313
$ws2.Range("I4").Value = 114.7833325408463
$ws2.Range("F5").Value = 1.865916773012622
$ws2.Range("G5").Value = 45.12020132643337
$ws2.Range("H5").Value = 0.1367381729688636
$ws2.Range("I5").Value = 15.85810510146202

# --- Sheet 3: "Species std" (A1:H3) ---
$ws3 = $wb.Worksheets.Item("Species std")
$ws3.Range("E2").Value = 0.1397954707334627
$ws3.Range("F2").Value = 2.39453723239675
$ws3.Range("G2").Value = 0.02841795706590938
$ws3.Range("H2").Value = 619.694887831065
$ws3.Range("E3").Value = 2.273331134430574
$ws3.Range("F3").Value = 42.36182888418334
$ws3.Range("G3").Value = 0.06775403755896021
$ws3.Range("H3").Value = 81.68345038146514

# --- Sheet 4: "Species bioshaker std" (A1:I5) ---
$ws4 = $wb.Worksheets.Item("Species bioshaker std")
$ws4.Range("F2").Value = 0.04286242116351226
$ws4.Range("G2").Value = 0.1645102921142492
$ws4.Range("H2").Value = 0.009326732077518744
$ws4.Range("I2").Value = 789.8194392222913
$ws4.Range("F3").Value = 0.1020158181379437
$ws4.Range("G3").Value = 1.47932918866399
$ws4.Range("H3").Value = 0.02149679771825231
$ws4.Range("I3").Value = 296.6903433996373
$ws4.Range("F4").Value = 0.107722811759297
$ws4.Range("G4").Value = 2.264400101539799
$ws4.Range("H4").Value = 0.008169488179894364
$ws4.Range("I4").Value = 87.95864345558805
$ws4.Range("F5").Value = 3.26645675810937
$ws4.Range("G5").Value = 61.12017615726971
$ws4.Range("H5").Value = 0.0734476398361644
$ws4.Range("I5").Value = 13.46880641448709
